# "Generate Report for Handoff" - refresh the localization-status report
# with the new handoff id (0f0e4593-... -> 8846ec48-...) and updated
# timestamps, exactly as the CI job that produces this workbook would.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item(1)   # "Overview"
$wsZhCn     = $wb.Worksheets.Item(2)   # "zh-cn"
$wsDeDe     = $wb.Worksheets.Item(3)   # "de-de"

$oldId = "0f0e4593-53fb-44f0-97c9-8f95d138bfc9"
$newId = "8846ec48-452a-42ae-9718-1bceb151b7dd"

$oldZhXlf = "$oldId.fb203cf0a6b4a81ca2d3f2ed8bdd4d427176c98e.zh-cn.xlf"
$newZhXlf = "$newId.edd285a4b9368bdddfab3071201006d8759cb230.zh-cn.xlf"

$oldDeXlf = "$oldId.fb203cf0a6b4a81ca2d3f2ed8bdd4d427176c98e.de-de.xlf"
$newDeXlf = "$newId.edd285a4b9368bdddfab3071201006d8759cb230.de-de.xlf"

$oldHandoffDate = "2016-08-25 10:59:45"
$newHandoffDate = "2016-08-25 11:00:18"

$oldZhDate = "2016-08-25 10:59:40"
$newZhDate = "2016-08-25 10:59:58"

$hyperlinkAddress = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e85c087944adcc3a9c711e24e84a1123cb3da29b/e2e/$oldId.md"

# --- Plain text / date cells: update every occurrence so the shared
# string for each unique value is rewritten consistently everywhere it
# is used. ---

$wsOverview.Range("A2").Value = "$newId.md"
$wsZhCn.Range("A2").Value = "$newId.md"
$wsDeDe.Range("A2").Value = "$newId.md"

$wsOverview.Range("G2").Value = $newHandoffDate
$wsDeDe.Range("H2").Value = $newHandoffDate

$wsZhCn.Range("G2").Value = $newZhXlf
$wsZhCn.Range("H2").Value = $newZhDate

$wsDeDe.Range("G2").Value = $newDeXlf

# --- Hyperlink cells: the displayed text needs to change too. Re-create
# each hyperlink (same target address, same relationship) so the
# "display" text is refreshed without leaving a stale duplicate link
# behind. ---

$wsOverview.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("B2"), $hyperlinkAddress, "", "", "e2e\$newId.md")

$wsZhCn.Hyperlinks.Delete()
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A2"), $hyperlinkAddress, "", "", "$newId.md")

$wsDeDe.Hyperlinks.Delete()
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A2"), $hyperlinkAddress, "", "", "$newId.md")
